# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the 9eb30ea2-... handback row (row 7) on both the zh-cn and
# de-de sheets, reflecting a regenerated handback report with later
# timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D7").Value = "2016-03-03 10:27:26"
$wsZhCn.Range("G7").Value = "2016-03-03 10:28:23"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D7").Value = "2016-03-03 10:27:38"
$wsDeDe.Range("G7").Value = "2016-03-03 10:28:46"
